{"js": "// 1) Merge the split \"on-napkins\" run back into a single run by\n//    re-typing the sentence that was previously split across three\n//    runs (\"...on-napki\" + \"n\" + \"s' design...\"). Content is unchanged,\n//    only the run that carries it is normalised back to one run.\nconst mergeTarget =\n  \"on-napkins\\u2019 design to one that can (and should be) specified denotatively with denotational design.\";\nconst mergeResults = context.document.body.search(mergeTarget, { matchCase: true });\nmergeResults.load(\"items,length\");\nawait context.sync();\nif (mergeResults.items.length > 0) {\n  mergeResults.items[0].insertText(mergeTarget, \"Replace\");\n  await context.sync();\n}\n\n// 2) Swap the argument order of the `Property -> Relation` axioms to\n//    `Relation -> Property` to match the updated implementation\n//    signature (get/set/getAsStream/setAsStream).\nconst argSwapResults = context.document.body.search(\"Property -> Relation\", { matchCase: true });\nargSwapResults.load(\"items,length\");\nawait context.sync();\nfor (let i = 0; i < argSwapResults.items.length; i++) {\n  argSwapResults.items[i].insertText(\"Relation -> Property\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Re-type the sentence that was previously split across three runs\n#    (\"...on-napki\" + \"n\" + \"s' design...\") so it collapses back into a\n#    single run. The visible text does not change.\n$mergeTarget = \"on-napkins\" + [char]8217 + \" design to one that can (and should be) specified denotatively with denotational design.\"\n$findMerge = $d.Content.Find\n$findMerge.Text = $mergeTarget\n$findMerge.Replacement.Text = $mergeTarget\n$findMerge.Execute([ref]$findMerge.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$findMerge.Replacement.Text, 2) | Out-Null\n\n# 2) Swap the argument order of the `Property -> Relation` axioms to\n#    `Relation -> Property` (get/set/getAsStream/setAsStream) to match\n#    the updated implementation signature.\n$findSwap = $d.Content.Find\n$findSwap.Text = \"Property -> Relation\"\n$findSwap.Replacement.Text = \"Relation -> Property\"\n$findSwap.Execute([ref]$findSwap.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$findSwap.Replacement.Text, 2) | Out-Null\n"}
